# Raw and Clean Data from SSA for August 14th (2020-08-14, Excel serial 44057)
# Appends a new daily row (row 76) to out_vars / dates_dx / dates_sx / dates_deaths,
# extends the control_obs tracking table with a new date column (BX), and leaves
# control_obs as the active sheet/tab (matching the author's final selection).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# out_vars (raw daily totals)
# ---------------------------------------------------------------------------
$wsOut = $wb.Worksheets.Item("out_vars")
$wsOut.Activate()

$wsOut.Range("A75:J75").Copy()
$wsOut.Range("A76:J76").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$wsOut.Range("A76").Value = 44057
$wsOut.Range("B76").Value = 511369
$wsOut.Range("C76").Value = 559974
$wsOut.Range("D76").Value = 85509
$wsOut.Range("E76").Value = 55908
$wsOut.Range("F76").Value = 26.325999425072698
$wsOut.Range("G76").Value = 134623
$wsOut.Range("H76").Value = 10798
$wsOut.Range("I76").Value = 13022
$wsOut.Range("J76").Value = 1156852
[void]$wsOut.Range("A76").Select()

# ---------------------------------------------------------------------------
# dates_dx (confirmed cases by symptom-onset-to-diagnosis lag bucket)
# ---------------------------------------------------------------------------
$wsDx = $wb.Worksheets.Item("dates_dx")
$wsDx.Activate()

$wsDx.Range("A75:L75").Copy()
$wsDx.Range("A76:L76").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$wsDx.Range("A76").Value = 44057
$wsDx.Range("B76").Value = 0
$wsDx.Range("C76").Value = 1
$wsDx.Range("D76").Value = 0
$wsDx.Range("E76").Value = 0
$wsDx.Range("F76").Value = 1
$wsDx.Range("G76").Value = 0
$wsDx.Range("H76").Value = 0
$wsDx.Range("I76").Value = 0
$wsDx.Range("J76").Value = 0
$wsDx.Range("K76").Value = 0
$wsDx.Range("L76").Value = 4
[void]$wsDx.Range("D77").Select()

# ---------------------------------------------------------------------------
# dates_sx (confirmed cases by symptom-onset lag bucket)
# ---------------------------------------------------------------------------
$wsSx = $wb.Worksheets.Item("dates_sx")
$wsSx.Activate()

$wsSx.Range("A75:N75").Copy()
$wsSx.Range("A76:N76").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$wsSx.Range("A76").Value = 44057
$wsSx.Range("B76").Value = 0
$wsSx.Range("C76").Value = 1
$wsSx.Range("D76").Value = 0
$wsSx.Range("E76").Value = 0
$wsSx.Range("F76").Value = 0
$wsSx.Range("G76").Value = 0
$wsSx.Range("H76").Value = 1
$wsSx.Range("I76").Value = 0
$wsSx.Range("J76").Value = 0
$wsSx.Range("K76").Value = 1
$wsSx.Range("L76").Value = 0
$wsSx.Range("M76").Value = 0
$wsSx.Range("N76").Value = 0
[void]$wsSx.Range("O76").Select()

# ---------------------------------------------------------------------------
# dates_deaths (deaths by lag bucket)
# ---------------------------------------------------------------------------
$wsDeaths = $wb.Worksheets.Item("dates_deaths")
$wsDeaths.Activate()

$wsDeaths.Range("A75:J75").Copy()
$wsDeaths.Range("A76:J76").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$wsDeaths.Range("A76").Value = 44057
$wsDeaths.Range("B76").Value = 0
$wsDeaths.Range("C76").Value = 0
$wsDeaths.Range("D76").Value = 0
$wsDeaths.Range("E76").Value = 0
$wsDeaths.Range("F76").Value = 2
$wsDeaths.Range("G76").Value = 1
$wsDeaths.Range("H76").Value = 1
$wsDeaths.Range("I76").Value = 1
$wsDeaths.Range("J76").Value = 2
[void]$wsDeaths.Range("J76").Select()

# ---------------------------------------------------------------------------
# control_obs (running QA counters) - new date column BX (2020-08-14 / 44057)
# ---------------------------------------------------------------------------
$wsCtrl = $wb.Worksheets.Item("control_obs")
$wsCtrl.Activate()

# Normalize BV1:BW1 formatting (they had drifted from the header style) and
# extend the header with the new date in the same style.
$wsCtrl.Range("BU1").Copy()
$wsCtrl.Range("BV1:BX1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0
$wsCtrl.Range("BV1").Value = 44055
$wsCtrl.Range("BW1").Value = 44056
$wsCtrl.Range("BX1").Value = 44057

$wsCtrl.Range("BW2:BW8").Copy()
$wsCtrl.Range("BX2:BX8").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0
$wsCtrl.Range("BX2").Value = 5314
$wsCtrl.Range("BX3").Value = 5113
$wsCtrl.Range("BX4").Value = 5113
$wsCtrl.Range("BX5").Value = 5113
$wsCtrl.Range("BX6").Value = 5113
$wsCtrl.Range("BX7").Value = 4423
$wsCtrl.Range("BX8").Value = 7054

$wsCtrl.Range("BW10:BW16").Copy()
$wsCtrl.Range("BX10:BX16").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0
$wsCtrl.Range("BX10").Value = 215
$wsCtrl.Range("BX11").Value = 215
$wsCtrl.Range("BX12").Value = 215
$wsCtrl.Range("BX13").Value = 215
$wsCtrl.Range("BX14").Value = 215
$wsCtrl.Range("BX15").Value = 150
$wsCtrl.Range("BX16").Value = 227

$wsCtrl.Range("BX18").Value = 1234

$wsCtrl.Range("BX20").Formula = "=SUM(BX2:BX18)"

[void]$wsCtrl.Range("BW26").Select()

$wb.Save()
Write-Output "edit applied"
